$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-18 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-19 Thursday", 2) | Out-Null
$d.Content.Find.Execute("32+40=", $true, $false, $false, $false, $false, $true, 1, $false, "42+17=", 2) | Out-Null
$d.Content.Find.Execute("43-33=", $true, $false, $false, $false, $false, $true, 1, $false, "78-32=", 2) | Out-Null
$d.Content.Find.Execute("13+24=", $true, $false, $false, $false, $false, $true, 1, $false, "94-85=", 2) | Out-Null
$d.Content.Find.Execute("62+20=", $true, $false, $false, $false, $false, $true, 1, $false, "13+38=", 2) | Out-Null
$d.Content.Find.Execute("12-1=", $true, $false, $false, $false, $false, $true, 1, $false, "28-12=", 2) | Out-Null
$d.Content.Find.Execute("46-9=", $true, $false, $false, $false, $false, $true, 1, $false, "81-9=", 2) | Out-Null
$d.Content.Find.Execute("97-54=", $true, $false, $false, $false, $false, $true, 1, $false, "55-16=", 2) | Out-Null
$d.Content.Find.Execute("84+7=", $true, $false, $false, $false, $false, $true, 1, $false, "72-62=", 2) | Out-Null
$d.Content.Find.Execute("79-46=", $true, $false, $false, $false, $false, $true, 1, $false, "27+57=", 2) | Out-Null
$d.Content.Find.Execute("4+57=", $true, $false, $false, $false, $false, $true, 1, $false, "77-16=", 2) | Out-Null
$d.Content.Find.Execute("35-0=", $true, $false, $false, $false, $false, $true, 1, $false, "96-73=", 2) | Out-Null
$d.Content.Find.Execute("8+22=", $true, $false, $false, $false, $false, $true, 1, $false, "3+9=", 2) | Out-Null
$d.Content.Find.Execute("67-62=", $true, $false, $false, $false, $false, $true, 1, $false, "73-42=", 2) | Out-Null
$d.Content.Find.Execute("85-53=", $true, $false, $false, $false, $false, $true, 1, $false, "23-9=", 2) | Out-Null
$d.Content.Find.Execute("72-44=", $true, $false, $false, $false, $false, $true, 1, $false, "84-9=", 2) | Out-Null
$d.Content.Find.Execute("76+9=", $true, $false, $false, $false, $false, $true, 1, $false, "87-70=", 2) | Out-Null
$d.Content.Find.Execute("48+35=", $true, $false, $false, $false, $false, $true, 1, $false, "50-44=", 2) | Out-Null
$d.Content.Find.Execute("77-53=", $true, $false, $false, $false, $false, $true, 1, $false, "73-31=", 2) | Out-Null
$d.Content.Find.Execute("63-53=", $true, $false, $false, $false, $false, $true, 1, $false, "60-42=", 2) | Out-Null
$d.Content.Find.Execute("35-23=", $true, $false, $false, $false, $false, $true, 1, $false, "18+70=", 2) | Out-Null
$d.Content.Find.Execute("37+29=", $true, $false, $false, $false, $false, $true, 1, $false, "62-26=", 2) | Out-Null
$d.Content.Find.Execute("56-2=", $true, $false, $false, $false, $false, $true, 1, $false, "59-35=", 2) | Out-Null
$d.Content.Find.Execute("38+59=", $true, $false, $false, $false, $false, $true, 1, $false, "72-51=", 2) | Out-Null
$d.Content.Find.Execute("97-66=", $true, $false, $false, $false, $false, $true, 1, $false, "86+7=", 2) | Out-Null
$d.Content.Find.Execute("40+35=", $true, $false, $false, $false, $false, $true, 1, $false, "66-39=", 2) | Out-Null
$d.Content.Find.Execute("65-62=", $true, $false, $false, $false, $false, $true, 1, $false, "41+13=", 2) | Out-Null
$d.Content.Find.Execute("67-42=", $true, $false, $false, $false, $false, $true, 1, $false, "10+6=", 2) | Out-Null
$d.Content.Find.Execute("12+36=", $true, $false, $false, $false, $false, $true, 1, $false, "41+2=", 2) | Out-Null
$d.Content.Find.Execute("34+62=", $true, $false, $false, $false, $false, $true, 1, $false, "54-10=", 2) | Out-Null
$d.Content.Find.Execute("30+60=", $true, $false, $false, $false, $false, $true, 1, $false, "91-0=", 2) | Out-Null
$d.Content.Find.Execute("79-22=", $true, $false, $false, $false, $false, $true, 1, $false, "16+79=", 2) | Out-Null
$d.Content.Find.Execute("50-25=", $true, $false, $false, $false, $false, $true, 1, $false, "30+66=", 2) | Out-Null
$d.Content.Find.Execute("24+20=", $true, $false, $false, $false, $false, $true, 1, $false, "80+5=", 2) | Out-Null
$d.Content.Find.Execute("54-19=", $true, $false, $false, $false, $false, $true, 1, $false, "79+19=", 2) | Out-Null
$d.Content.Find.Execute("46+25=", $true, $false, $false, $false, $false, $true, 1, $false, "73-56=", 2) | Out-Null
$d.Content.Find.Execute("33+1=", $true, $false, $false, $false, $false, $true, 1, $false, "54+14=", 2) | Out-Null
$d.Content.Find.Execute("53-37=", $true, $false, $false, $false, $false, $true, 1, $false, "73-27=", 2) | Out-Null
$d.Content.Find.Execute("78-24=", $true, $false, $false, $false, $false, $true, 1, $false, "60-23=", 2) | Out-Null
$d.Content.Find.Execute("72-63=", $true, $false, $false, $false, $false, $true, 1, $false, "90-42=", 2) | Out-Null
$d.Content.Find.Execute("47+22=", $true, $false, $false, $false, $false, $true, 1, $false, "78-58=", 2) | Out-Null
$d.Content.Find.Execute("44+3=", $true, $false, $false, $false, $false, $true, 1, $false, "64-10=", 2) | Out-Null
$d.Content.Find.Execute("99-2=", $true, $false, $false, $false, $false, $true, 1, $false, "49-19=", 2) | Out-Null
$d.Content.Find.Execute("25+68=", $true, $false, $false, $false, $false, $true, 1, $false, "3+35=", 2) | Out-Null
$d.Content.Find.Execute("18+3=", $true, $false, $false, $false, $false, $true, 1, $false, "77-16=", 2) | Out-Null
$d.Content.Find.Execute("81-41=", $true, $false, $false, $false, $false, $true, 1, $false, "76-3=", 2) | Out-Null
$d.Content.Find.Execute("38-2=", $true, $false, $false, $false, $false, $true, 1, $false, "57-21=", 2) | Out-Null
$d.Content.Find.Execute("16+31=", $true, $false, $false, $false, $false, $true, 1, $false, "21+41=", 2) | Out-Null
$d.Content.Find.Execute("89-78=", $true, $false, $false, $false, $false, $true, 1, $false, "21+37=", 2) | Out-Null
$d.Content.Find.Execute("99-90=", $true, $false, $false, $false, $false, $true, 1, $false, "28+66=", 2) | Out-Null
$d.Content.Find.Execute("89-31=", $true, $false, $false, $false, $false, $true, 1, $false, "51+47=", 2) | Out-Null
$d.Content.Find.Execute("46+14=", $true, $false, $false, $false, $false, $true, 1, $false, "76-22=", 2) | Out-Null
$d.Content.Find.Execute("8-2=", $true, $false, $false, $false, $false, $true, 1, $false, "14+3=", 2) | Out-Null
$d.Content.Find.Execute("20+14=", $true, $false, $false, $false, $false, $true, 1, $false, "71-37=", 2) | Out-Null
$d.Content.Find.Execute("58-4=", $true, $false, $false, $false, $false, $true, 1, $false, "92-52=", 2) | Out-Null
$d.Content.Find.Execute("56+24=", $true, $false, $false, $false, $false, $true, 1, $false, "76-50=", 2) | Out-Null
$d.Content.Find.Execute("85+6=", $true, $false, $false, $false, $false, $true, 1, $false, "86+12=", 2) | Out-Null
$d.Content.Find.Execute("18-13=", $true, $false, $false, $false, $false, $true, 1, $false, "84-29=", 2) | Out-Null
$d.Content.Find.Execute("23+29=", $true, $false, $false, $false, $false, $true, 1, $false, "12+14=", 2) | Out-Null
$d.Content.Find.Execute("88-7=", $true, $false, $false, $false, $false, $true, 1, $false, "9+82=", 2) | Out-Null
$d.Content.Find.Execute("40+38=", $true, $false, $false, $false, $false, $true, 1, $false, "39-3=", 2) | Out-Null
$d.Content.Find.Execute("84-33=", $true, $false, $false, $false, $false, $true, 1, $false, "92-46=", 2) | Out-Null
$d.Content.Find.Execute("59-5=", $true, $false, $false, $false, $false, $true, 1, $false, "20+12=", 2) | Out-Null
$d.Content.Find.Execute("5+48=", $true, $false, $false, $false, $false, $true, 1, $false, "90-74=", 2) | Out-Null
$d.Content.Find.Execute("27-26=", $true, $false, $false, $false, $false, $true, 1, $false, "50-33=", 2) | Out-Null
$d.Content.Find.Execute("55+22=", $true, $false, $false, $false, $false, $true, 1, $false, "28+38=", 2) | Out-Null
$d.Content.Find.Execute("20+49=", $true, $false, $false, $false, $false, $true, 1, $false, "7+4=", 2) | Out-Null
$d.Content.Find.Execute("35+45=", $true, $false, $false, $false, $false, $true, 1, $false, "2+27=", 2) | Out-Null
$d.Content.Find.Execute("71+6=", $true, $false, $false, $false, $false, $true, 1, $false, "6+89=", 2) | Out-Null
$d.Content.Find.Execute("69-50=", $true, $false, $false, $false, $false, $true, 1, $false, "53-7=", 2) | Out-Null
$d.Content.Find.Execute("30-4=", $true, $false, $false, $false, $false, $true, 1, $false, "31+48=", 2) | Out-Null
$d.Content.Find.Execute("26+32=", $true, $false, $false, $false, $false, $true, 1, $false, "76+1=", 2) | Out-Null
$d.Content.Find.Execute("26+46=", $true, $false, $false, $false, $false, $true, 1, $false, "60-15=", 2) | Out-Null
$d.Content.Find.Execute("7+8=", $true, $false, $false, $false, $false, $true, 1, $false, "0+41=", 2) | Out-Null
$d.Content.Find.Execute("9+29=", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=", 2) | Out-Null
$d.Content.Find.Execute("56+36=", $true, $false, $false, $false, $false, $true, 1, $false, "71-58=", 2) | Out-Null
$d.Content.Find.Execute("4+47=", $true, $false, $false, $false, $false, $true, 1, $false, "72-24=", 2) | Out-Null
$d.Content.Find.Execute("25+43=", $true, $false, $false, $false, $false, $true, 1, $false, "99-42=", 2) | Out-Null
$d.Content.Find.Execute("18+57=", $true, $false, $false, $false, $false, $true, 1, $false, "7+70=", 2) | Out-Null
$d.Content.Find.Execute("93-84=", $true, $false, $false, $false, $false, $true, 1, $false, "88-30=", 2) | Out-Null
$d.Content.Find.Execute("13-3=", $true, $false, $false, $false, $false, $true, 1, $false, "64-13=", 2) | Out-Null
$d.Content.Find.Execute("82-37=", $true, $false, $false, $false, $false, $true, 1, $false, "15+66=", 2) | Out-Null
$d.Content.Find.Execute("38-9=", $true, $false, $false, $false, $false, $true, 1, $false, "70+10=", 2) | Out-Null
$d.Content.Find.Execute("9-9=", $true, $false, $false, $false, $false, $true, 1, $false, "40-2=", 2) | Out-Null
$d.Content.Find.Execute("99-66=", $true, $false, $false, $false, $false, $true, 1, $false, "39+36=", 2) | Out-Null
$d.Content.Find.Execute("33+18=", $true, $false, $false, $false, $false, $true, 1, $false, "6+59=", 2) | Out-Null
$d.Content.Find.Execute("50+41=", $true, $false, $false, $false, $false, $true, 1, $false, "78+11=", 2) | Out-Null
$d.Content.Find.Execute("10-8=", $true, $false, $false, $false, $false, $true, 1, $false, "58-33=", 2) | Out-Null
$d.Content.Find.Execute("23+32=", $true, $false, $false, $false, $false, $true, 1, $false, "31+4=", 2) | Out-Null
$d.Content.Find.Execute("94-1=", $true, $false, $false, $false, $false, $true, 1, $false, "21+70=", 2) | Out-Null
$d.Content.Find.Execute("71+28=", $true, $false, $false, $false, $false, $true, 1, $false, "69-20=", 2) | Out-Null
$d.Content.Find.Execute("54-6=", $true, $false, $false, $false, $false, $true, 1, $false, "60-45=", 2) | Out-Null
$d.Content.Find.Execute("61+28=", $true, $false, $false, $false, $false, $true, 1, $false, "0+60=", 2) | Out-Null
$d.Content.Find.Execute("67+12=", $true, $false, $false, $false, $false, $true, 1, $false, "21+20=", 2) | Out-Null
$d.Content.Find.Execute("20+73=", $true, $false, $false, $false, $false, $true, 1, $false, "80-0=", 2) | Out-Null
$d.Content.Find.Execute("34+24=", $true, $false, $false, $false, $false, $true, 1, $false, "19+76=", 2) | Out-Null
$d.Content.Find.Execute("95-75=", $true, $false, $false, $false, $false, $true, 1, $false, "87-42=", 2) | Out-Null
$d.Content.Find.Execute("26+59=", $true, $false, $false, $false, $false, $true, 1, $false, "27+53=", 2) | Out-Null
$d.Content.Find.Execute("5+31=", $true, $false, $false, $false, $false, $true, 1, $false, "30+5=", 2) | Out-Null
$d.Content.Find.Execute("2+35=", $true, $false, $false, $false, $false, $true, 1, $false, "43+1=", 2) | Out-Null
$d.Content.Find.Execute("76-20=", $true, $false, $false, $false, $false, $true, 1, $false, "93-91=", 2) | Out-Null

Write-Output "Replacements applied"